$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = 117
$ws.Range("D14").Value = 143
$ws.Range("D15").Value = 43

$ws.Range("E13:E15").Select()
